$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing year column (S) into the new
# column (T) before writing the 2023 figures, so the new column inherits
# the same borders/fonts/number formats as the rest of the table.
$ws.Range("S4:S8").Copy()
$ws.Range("T4:T8").PasteSpecial(-4122)

$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 43.1
$ws.Range("T6").Value = 19.7
$ws.Range("T7").Value = 7.8
$ws.Range("T8").Value = 15.6

# Row height tweaks that came with the new column being added.
$ws.Rows.Item(1).RowHeight = 57
$ws.Rows.Item(4).RowHeight = 16.5

# Clear the leftover interactive selection (was sitting on Y14, well past
# the used range) and park it back on A1.
$ws.Range("A1").Select()
